# correção nos dados e inicio da analise PNAD 2009
# Remove the category-header / footnote rows that have no data, so the
# remaining category rows (homens, mulheres, branca, ...) shift up and
# keep their original figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(36, 35, 29, 27, 19, 13, 8, 5)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
